# Apply updated cryptocurrency price/volume figures to the symbol list.
# Each target cell stores its value as literal text (matching the source
# data's inline-string cells), so we prefix the new value with a leading
# apostrophe to force Excel to keep it as text rather than coercing it to
# a Number/Percentage (which would lose exact decimal formatting, e.g.
# trailing zeros, and introduce floating point rounding noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'277.93"
$ws.Range("E2").Value = "'6.08%"

# Row 3
$ws.Range("E3").Value = "'2.62%"

# Row 4
$ws.Range("D4").Value = "'4.801"

# Row 5
$ws.Range("D5").Value = "'0.06339"
$ws.Range("E5").Value = "'4.29%"

# Row 6
$ws.Range("D6").Value = "'6.939"
$ws.Range("E6").Value = "'3.52%"

# Row 7
$ws.Range("D7").Value = "'3.359"
$ws.Range("E7").Value = "'5.97%"

# Row 8
$ws.Range("D8").Value = "'0.8823"
$ws.Range("E8").Value = "'3.68%"

# Row 9
$ws.Range("D9").Value = "'0.9476"
$ws.Range("E9").Value = "'4.33%"

# Row 10
$ws.Range("D10").Value = "'0.1479"
$ws.Range("E10").Value = "'5.54%"

# Row 11
$ws.Range("D11").Value = "'0.05269"
$ws.Range("E11").Value = "'3.44%"

# Row 12
$ws.Range("D12").Value = "'0.07332"
$ws.Range("E12").Value = "'3.38%"

# Row 13
$ws.Range("D13").Value = "'0.03127"
$ws.Range("E13").Value = "'0.30%"

# Row 14
$ws.Range("D14").Value = "'0.09066"
$ws.Range("E14").Value = "'0.25%"

# Row 15
$ws.Range("D15").Value = "'0.001553"
$ws.Range("E15").Value = "'1.35%"

# Row 16
$ws.Range("D16").Value = "'0.0006262"
$ws.Range("E16").Value = "'1.66%"

# Row 17
$ws.Range("D17").Value = "'0.005774"
$ws.Range("E17").Value = "'-3.48%"

# Row 18
$ws.Range("D18").Value = "'3.463"
$ws.Range("E18").Value = "'0.38%"

# Row 19
$ws.Range("D19").Value = "'2.287"
$ws.Range("E19").Value = "'6.57%"

# Row 20
$ws.Range("D20").Value = "'0.3096"
$ws.Range("E20").Value = "'0.77%"

# Row 21
$ws.Range("D21").Value = "'0.1338"
$ws.Range("E21").Value = "'4.46%"

# Row 22
$ws.Range("D22").Value = "'3.866"
$ws.Range("E22").Value = "'-6.17%"

# Row 23
$ws.Range("D23").Value = "'0.04320"
$ws.Range("E23").Value = "'2.17%"

# Row 24
$ws.Range("D24").Value = "'0.001175"
$ws.Range("E24").Value = "'-0.33%"

# Row 25
$ws.Range("D25").Value = "'0.003587"
$ws.Range("E25").Value = "'-11.64%"

# Row 27
$ws.Range("D27").Value = "'0.0001688"
$ws.Range("E27").Value = "'-12.91%"

# Row 40
$ws.Range("D40").Value = "'0.04075"
$ws.Range("E40").Value = "'3.00%"

# Row 41
$ws.Range("D41").Value = "'0.006647"
$ws.Range("E41").Value = "'58.78%"

# Row 42
$ws.Range("D42").Value = "'0.1166"
$ws.Range("E42").Value = "'4.79%"

# Row 43
$ws.Range("D43").Value = "'0.002296"
$ws.Range("E43").Value = "'14.19%"

# Row 44
$ws.Range("D44").Value = "'0.01244"
$ws.Range("E44").Value = "'-3.94%"

# Row 45
$ws.Range("D45").Value = "'0.00005211"
$ws.Range("E45").Value = "'1.74%"

# Row 46
$ws.Range("E46").Value = "'-0.06%"

# Row 47
$ws.Range("D47").Value = "'2.379"
$ws.Range("E47").Value = "'820.54%"

# Row 48
$ws.Range("D48").Value = "'0.02247"
$ws.Range("E48").Value = "'5.92%"

# Row 49
$ws.Range("E49").Value = "'-0.06%"

# Row 50
$ws.Range("E50").Value = "'-0.13%"
